$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Gas6"
$ws.Range("C2").Value = "Tyro3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 12.28498633333333
$ws.Range("H2").Value = 36.854959
$ws.Range("I2").Value = 0.1279589698403688
$ws.Range("J2").Value = 0.1279589698403688
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.07610666666666667
$ws.Range("N2").Value = 0.22832
$ws.Range("O2").Value = 0.04153709822108353
$ws.Range("P2").Value = 0.04153709822108354
$ws.Range("Q2").Value = 0.9349693598755556
$ws.Range("R2").Value = 8.414724238880002
$ws.Range("S2").Value = 0.005315044298528066
$ws.Range("T2").Value = 0.005315044298528066

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Gas6"
$ws.Range("C3").Value = "Tyro3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 12.28498633333333
$ws.Range("H3").Value = 36.854959
$ws.Range("I3").Value = 0.1279589698403688
$ws.Range("J3").Value = 0.1279589698403688
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.639501
$ws.Range("N3").Value = 4.918502999999999
$ws.Range("O3").Value = 0.8947982752789682
$ws.Range("P3").Value = 0.8947982752789683
$ws.Range("Q3").Value = 20.14124737848633
$ws.Range("R3").Value = 181.271226406377
$ws.Range("S3").Value = 0.1144974655196355
$ws.Range("T3").Value = 0.1144974655196355

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Gas6"
$ws.Range("C4").Value = "Tyro3"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 12.28498633333333
$ws.Range("H4").Value = 36.854959
$ws.Range("I4").Value = 0.1279589698403688
$ws.Range("J4").Value = 0.1279589698403688
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.11665
$ws.Range("N4").Value = 0.34995
$ws.Range("O4").Value = 0.06366462649994824
$ws.Range("P4").Value = 0.06366462649994825
$ws.Range("Q4").Value = 1.433043655783333
$ws.Range("R4").Value = 12.89739290205
$ws.Range("S4").Value = 0.008146460022205223
$ws.Range("T4").Value = 0.008146460022205223

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Gas6"
$ws.Range("C5").Value = "Tyro3"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 65.605124
$ws.Range("H5").Value = 196.815372
$ws.Range("I5").Value = 0.6833352399026945
$ws.Range("J5").Value = 0.6833352399026945
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.07610666666666667
$ws.Range("N5").Value = 0.22832
$ws.Range("O5").Value = 0.04153709822108353
$ws.Range("P5").Value = 0.04153709822108354
$ws.Range("Q5").Value = 4.992987303893334
$ws.Range("R5").Value = 44.93688573504
$ws.Range("S5").Value = 0.0283837629777659
$ws.Range("T5").Value = 0.02838376297776591

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Gas6"
$ws.Range("C6").Value = "Tyro3"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 65.605124
$ws.Range("H6").Value = 196.815372
$ws.Range("I6").Value = 0.6833352399026945
$ws.Range("J6").Value = 0.6833352399026945
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.639501
$ws.Range("N6").Value = 4.918502999999999
$ws.Range("O6").Value = 0.8947982752789682
$ws.Range("P6").Value = 0.8947982752789683
$ws.Range("Q6").Value = 107.559666403124
$ws.Range("R6").Value = 968.0369976281158
$ws.Range("S6").Value = 0.611447194102271
$ws.Range("T6").Value = 0.6114471941022711

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Gas6"
$ws.Range("C7").Value = "Tyro3"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 65.605124
$ws.Range("H7").Value = 196.815372
$ws.Range("I7").Value = 0.6833352399026945
$ws.Range("J7").Value = 0.6833352399026945
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.11665
$ws.Range("N7").Value = 0.34995
$ws.Range("O7").Value = 0.06366462649994824
$ws.Range("P7").Value = 0.06366462649994825
$ws.Range("Q7").Value = 7.6528377146
$ws.Range("R7").Value = 68.8755394314
$ws.Range("S7").Value = 0.04350428282265757
$ws.Range("T7").Value = 0.04350428282265758

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Gas6"
$ws.Range("C8").Value = "Tyro3"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 18.11712033333333
$ws.Range("H8").Value = 54.351361
$ws.Range("I8").Value = 0.1887057902569366
$ws.Range("J8").Value = 0.1887057902569366
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.07610666666666667
$ws.Range("N8").Value = 0.22832
$ws.Range("O8").Value = 0.04153709822108353
$ws.Range("P8").Value = 0.04153709822108354
$ws.Range("Q8").Value = 1.378833638168889
$ws.Range("R8").Value = 12.40950274352
$ws.Range("S8").Value = 0.007838290944789564
$ws.Range("T8").Value = 0.007838290944789566

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Gas6"
$ws.Range("C9").Value = "Tyro3"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 18.11712033333333
$ws.Range("H9").Value = 54.351361
$ws.Range("I9").Value = 0.1887057902569366
$ws.Range("J9").Value = 0.1887057902569366
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.639501
$ws.Range("N9").Value = 4.918502999999999
$ws.Range("O9").Value = 0.8947982752789682
$ws.Range("P9").Value = 0.8947982752789683
$ws.Range("Q9").Value = 29.70303690362033
$ws.Range("R9").Value = 267.327332132583
$ws.Range("S9").Value = 0.1688536156570616
$ws.Range("T9").Value = 0.1688536156570616

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Gas6"
$ws.Range("C10").Value = "Tyro3"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 18.11712033333333
$ws.Range("H10").Value = 54.351361
$ws.Range("I10").Value = 0.1887057902569366
$ws.Range("J10").Value = 0.1887057902569366
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.11665
$ws.Range("N10").Value = 0.34995
$ws.Range("O10").Value = 0.06366462649994824
$ws.Range("P10").Value = 0.06366462649994825
$ws.Range("Q10").Value = 2.113362086883333
$ws.Range("R10").Value = 19.02025878195
$ws.Range("S10").Value = 0.01201388365508544
$ws.Range("T10").Value = 0.01201388365508544
